# Reproduce the "Add files via upload" commit:
#   - rename the existing sheet "データ2" -> "data"
#   - add a new sheet "setting" right after it, containing 6 time-of-day
#     values in column A (A1:A6), formatted as h:mm
#   - the workbook-level _xlnm._FilterDatabase defined name follows the
#     sheet rename automatically, since it is a live reference

$wb = $excel.ActiveWorkbook

# --- rename the first (and only) worksheet -------------------------------
$wsData = $wb.Worksheets.Item(1)
$wsData.Name = "data"

# --- insert the new "setting" sheet right after "data" -------------------
$wsSetting = $wb.Worksheets.Add($null, $wsData)
$wsSetting.Name = "setting"

# --- fill A1:A6 with the time values, formatted as h:mm ------------------
# (stored as Excel serial-date fractions of a day, same as the source file)
$times = @(
    0.35416666666666669,  # 08:30
    0.5,                   # 12:00
    0.53472222222222221,  # 12:50
    0.64583333333333337,  # 15:30
    0.65277777777777779,  # 15:40
    0.72222222222222221   # 17:20
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $cell = $wsSetting.Cells.Item($i + 1, 1)
    $cell.Value = $times[$i]
    $cell.NumberFormat = "h:mm"
}

# --- restore the view state: whole-column selection on "setting", -------
# --- and keep "data" as the active/selected sheet -------------------------
$wsSetting.Range("A1:XFD1048576").Select()
$wsData.Select()
$wsData.Range("A1").Select()
